# Add new columns I and J ("I0" and "IF") to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - mirror the style used by the other header cells (B1:H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I ("I0") and J ("IF"), rows 2-14
$data = @{
    2  = @(7, 8)
    3  = @(4, 6)
    4  = @(3, 5)
    5  = @(8, 9)
    6  = @(9, 9)
    7  = @(6, 6)
    8  = @(9, 9)
    9  = @(5, 5)
    10 = @(5, 6)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(1, 1)
    14 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
